$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("O4").Value = 2023
Write-Output $ws.Range("O4").Value
